$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 3. This shifts the existing rows
# 3..77 down to 4..78 (preserving all of their data/formatting), which is
# exactly the "everything moved down one" pattern seen in the diff, and
# also grows the sheet dimension from A1:R77 to A1:R78 automatically.
$ws.Rows(3).Insert()

# Populate the newly-inserted row 3 with the new weekly record.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C3").Value = 'Ñuble'
$ws.Range("D3").Value = 44860
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112013
$ws.Range("G3").Value = 'Alcachofa'
$ws.Range("H3").Value = 'Española'
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 10500
$ws.Range("N3").Value = '$/caja 30 unidades'
$ws.Range("O3").Value = 'Provincia de Limarí'
$ws.Range("P3").Value = 350
$ws.Range("Q3").Value = 30
$ws.Range("R3").Value = 'Hortaliza'
